$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.708.19'
$ws.Range('E2').Value = '  +0.37%  '

$ws.Range('D3').Value = '2.511.82'
$ws.Range('E3').Value = '  +0.31%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.74'
$ws.Range('E5').Value = '  -0.82%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.08'
$ws.Range('E6').Value = '  -0.06%  '

$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('E8').Value = '  +0.18%  '

$ws.Range('D9').Value = '2.510.77'
$ws.Range('E9').Value = '  +0.18%  '

$ws.Range('E10').Value = '  +1.72%  '

$ws.Range('E11').Value = '  -0.57%  '

$ws.Range('E12').Value = '  +6.47%  '

$ws.Range('E13').Value = '  +1.84%  '

$ws.Range('D14').Value = '2.973.35'
$ws.Range('E14').Value = '  +0.24%  '

$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '69.637.02'
$ws.Range('E15').Value = '  +0.27%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000178'
$ws.Range('E16').Value = '  +1.50%  '

$ws.Range('E17').Value = '  +0.01%  '

$ws.Range('D18').Value = '2.523.33'
$ws.Range('E18').Value = '  +0.80%  '

$ws.Range('E19').Value = '  -0.59%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.60'
$ws.Range('E20').Value = '  -1.93%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '352.07'
$ws.Range('E21').Value = '  +0.06%  '

$ws.Range('E22').Value = '  -0.67%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.98'
$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.12'
$ws.Range('E24').Value = '  +2.89%  '

$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.19%  '

$ws.Range('E26').Value = '  -1.45%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.86'
$ws.Range('E27').Value = '  -2.18%  '

$ws.Range('D28').Value = '2.666.01'

$ws.Range('E29').Value = '  +0.63%  '

$ws.Range('E30').Value = '  -1.07%  '

$ws.Range('E31').Value = '  +0.01%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '462.63'
$ws.Range('E32').Value = '  -4.37%  '

$ws.Range('E33').Value = '  -5.38%  '

$ws.Range('E34').Value = '  -1.09%  '

$ws.Range('E35').Value = '  -0.04%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.02'
$ws.Range('E36').Value = '  +4.71%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.117'
$ws.Range('E37').Value = '  +1.11%  '

$ws.Range('E38').Value = '  +1.02%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.54'
$ws.Range('E39').Value = '  -0.23%  '

$ws.Range('E40').Value = '  +0.02%  '

$ws.Range('E41').Value = '  -0.23%  '

$ws.Range('E42').Value = '  -1.72%  '

$ws.Range('E43').Value = '  -1.49%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.22'
$ws.Range('E44').Value = '  -0.01%  '

$ws.Range('E45').Value = '  -5.13%  '

$ws.Range('E46').Value = '  -6.59%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.58'
$ws.Range('E47').Value = '  -0.73%  '

$ws.Range('E48').Value = '  -1.60%  '

$ws.Range('E49').Value = '  -1.60%  '

$ws.Range('E50').Value = '  +0.65%  '

$ws.Range('E51').Value = '  +3.08%  '
